$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("alpha_school_climate")
$ws.Range("B2").Value = 0.93811986667109926
$ws.Range("B5").Value = 0.52021457595164078
$ws.Range("C5").Value = 0.46088335413329407
$ws.Range("D5").Value = 0.44208830353793116
$ws.Range("E5").Value = 0.93771631321975168
$ws.Range("B6").Value = 0.59109636684383315
$ws.Range("C6").Value = 0.5386377017423325
$ws.Range("D6").Value = 0.43751734969827177
$ws.Range("E6").Value = 0.93662390329116285
$ws.Range("B7").Value = 0.47311259405578543
$ws.Range("C7").Value = 0.41138034118202949
$ws.Range("D7").Value = 0.44551708765766174
$ws.Range("E7").Value = 0.93852267399053968
$ws.Range("B8").Value = 0.54685124704182164
$ws.Range("C8").Value = 0.49043465758161364
$ws.Range("D8").Value = 0.44049910582768553
$ws.Range("E8").Value = 0.93733879463424963
$ws.Range("B9").Value = 0.53862495167377233
$ws.Range("C9").Value = 0.47094554250875525
$ws.Range("D9").Value = 0.44396287469296303
$ws.Range("E9").Value = 0.93815853523046633
$ws.Range("B10").Value = 0.78476474115159434
$ws.Range("C10").Value = 0.74717071260263157
$ws.Range("D10").Value = 0.42142270881698352
$ws.Range("E10").Value = 0.93261080169809207
$ws.Range("B11").Value = 0.8152789509224877
$ws.Range("C11").Value = 0.78282266201245565
$ws.Range("D11").Value = 0.41841427076689441
$ws.Range("E11").Value = 0.93183043203103233
$ws.Range("B12").Value = 0.76951665233257049
$ws.Range("C12").Value = 0.72910967655718306
$ws.Range("D12").Value = 0.42290307380571585
$ws.Range("E12").Value = 0.93299119636126671
$ws.Range("B13").Value = 0.84021421332046153
$ws.Range("C13").Value = 0.8086126296609738
$ws.Range("D13").Value = 0.41630576233609057
$ws.Range("E13").Value = 0.93127756733082889
$ws.Range("B14").Value = 0.7378263487386848
$ws.Range("C14").Value = 0.69203546312319986
$ws.Range("D14").Value = 0.42589317955654699
$ws.Range("E14").Value = 0.93375239972217805
$ws.Range("B15").Value = 0.76890380732505403
$ws.Range("C15").Value = 0.72755326672862342
$ws.Range("D15").Value = 0.42282378705797896
$ws.Range("E15").Value = 0.93297088250453564
$ws.Range("B16").Value = 0.82498666731825854
$ws.Range("C16").Value = 0.7943687624796536
$ws.Range("D16").Value = 0.41755781209995524
$ws.Range("E16").Value = 0.93160645726809987
$ws.Range("B17").Value = 0.81360892089385817
$ws.Range("C17").Value = 0.78035305157606061
$ws.Range("D17").Value = 0.41942866849568017
$ws.Range("E17").Value = 0.9320946649843429
$ws.Range("B18").Value = 0.26596340886531933
$ws.Range("C18").Value = 0.18794090417335127
$ws.Range("D18").Value = 0.46422552969655684
$ws.Range("E18").Value = 0.94273502981403401
$ws.Range("B19").Value = 0.65976071460772967
$ws.Range("C19").Value = 0.61335628927590491
$ws.Range("D19").Value = 0.43222507541203486
$ws.Range("E19").Value = 0.93533353035476741
$ws.Range("B20").Value = 0.61033157636876678
$ws.Range("C20").Value = 0.55488810277002099
$ws.Range("D20").Value = 0.43672983408746191
$ws.Range("E20").Value = 0.93643364671490026
$ws.Range("B21").Value = 0.74332002454423496
$ws.Range("C21").Value = 0.70281443410571387
$ws.Range("D21").Value = 0.42658145967339672
$ws.Range("E21").Value = 0.93392628118674414
$ws.Range("B22").Value = 0.75955882160231458
$ws.Range("C22").Value = 0.724805030149963
$ws.Range("D22").Value = 0.42606529162018753
$ws.Range("E22").Value = 0.93379592730747929
$ws.Range("B23").Value = 0.67642171967914766
$ws.Range("C23").Value = 0.62971163322889856
$ws.Range("D23").Value = 0.43162552827468886
$ws.Range("E23").Value = 0.93518557967801208
$ws.Range("B24").Value = 0.69513450010582534
$ws.Range("C24").Value = 0.64966957057701935
$ws.Range("D24").Value = 0.43047783149984775
$ws.Range("E24").Value = 0.93490134358069532

$ws = $wb.Worksheets.Item("alpha_teacher_quality")
$ws.Range("B2").Value = 0.89980548376827996
$ws.Range("B5").Value = 0.60765057404749367
$ws.Range("C5").Value = 0.54099130370298742
$ws.Range("D5").Value = 0.34711186131012206
$ws.Range("E5").Value = 0.89480876128231535
$ws.Range("B6").Value = 0.57890901691872954
$ws.Range("C6").Value = 0.50613589588786334
$ws.Range("D6").Value = 0.34931225267703742
$ws.Range("E6").Value = 0.89571783295340035
$ws.Range("B7").Value = 0.5704209036028508
$ws.Range("C7").Value = 0.50065923643525823
$ws.Range("D7").Value = 0.34918428976873922
$ws.Range("E7").Value = 0.89566522970794993
$ws.Range("B8").Value = 0.53503182891015633
$ws.Range("C8").Value = 0.45733791410152069
$ws.Range("D8").Value = 0.35260948418173588
$ws.Range("E8").Value = 0.89706219246221153
$ws.Range("B9").Value = 0.54932073727996222
$ws.Range("C9").Value = 0.47584282472940015
$ws.Range("D9").Value = 0.35127111453286131
$ws.Range("E9").Value = 0.89651906501406109
$ws.Range("B10").Value = 0.76272895271912244
$ws.Range("C10").Value = 0.7117936217212274
$ws.Range("D10").Value = 0.33088834881486806
$ws.Range("E10").Value = 0.88779554608734645
$ws.Range("B11").Value = 0.79083360931839053
$ws.Range("C11").Value = 0.74635828466019827
$ws.Range("D11").Value = 0.3283942465983597
$ws.Range("E11").Value = 0.88666629527227803
$ws.Range("B12").Value = 0.74898648072644181
$ws.Range("C12").Value = 0.69562621214810938
$ws.Range("D12").Value = 0.33254915084140568
$ws.Range("E12").Value = 0.88853967872568107
$ws.Range("B13").Value = 0.73765275982885248
$ws.Range("C13").Value = 0.68385891674021859
$ws.Range("D13").Value = 0.33440865992823343
$ws.Range("E13").Value = 0.88936553045657074
$ws.Range("B14").Value = 0.77405707165228466
$ws.Range("C14").Value = 0.72439603017699539
$ws.Range("D14").Value = 0.32982355817976178
$ws.Range("E14").Value = 0.8873151804133691
$ws.Range("B15").Value = 0.7321382045732433
$ws.Range("C15").Value = 0.6767540676300533
$ws.Range("D15").Value = 0.33478932057131988
$ws.Range("E15").Value = 0.88953364742748453
$ws.Range("B16").Value = 0.44061752009594934
$ws.Range("C16").Value = 0.34975492578928696
$ws.Range("D16").Value = 0.36286363371859487
$ws.Range("E16").Value = 0.90111114928563685
$ws.Range("B17").Value = 0.47620613915187499
$ws.Range("C17").Value = 0.38898181168285684
$ws.Range("D17").Value = 0.35988615307866706
$ws.Range("E17").Value = 0.89995551667692608
$ws.Range("B18").Value = 0.55467051167483483
$ws.Range("C18").Value = 0.47558276939538946
$ws.Range("D18").Value = 0.35284101886260411
$ws.Range("E18").Value = 0.89715580056323885
$ws.Range("B19").Value = 0.52646442951524608
$ws.Range("C19").Value = 0.44420698421811888
$ws.Range("D19").Value = 0.35552601201254852
$ws.Range("E19").Value = 0.89823382852846245
$ws.Range("B20").Value = 0.49350866726822967
$ws.Range("C20").Value = 0.40812797146854823
$ws.Range("D20").Value = 0.35838494836191814
$ws.Range("E20").Value = 0.89936672332607048
$ws.Range("B21").Value = 0.64001003189013472
$ws.Range("C21").Value = 0.57090761410195578
$ws.Range("D21").Value = 0.34507033591622038
$ws.Range("E21").Value = 0.8939566341337134

$ws = $wb.Worksheets.Item("alpha_student_support")
$ws.Range("B2").Value = 0.66211117647631301
$ws.Range("B5").Value = 0.72440688132534969
$ws.Range("C5").Value = 0.45125479042373434
$ws.Range("D5").Value = 0.31803270388592036
$ws.Range("E5").Value = 0.58316623964145087
$ws.Range("B6").Value = 0.64920938872831568
$ws.Range("C6").Value = 0.34990273077469564
$ws.Range("D6").Value = 0.38755529935631444
$ws.Range("E6").Value = 0.65498222979015974
$ws.Range("B7").Value = 0.76412230331561881
$ws.Range("C7").Value = 0.46341058296824783
$ws.Range("D7").Value = 0.29898614071388335
$ws.Range("E7").Value = 0.56131037600992029
$ws.Range("B8").Value = 0.74125341506934572
$ws.Range("C8").Value = 0.46009249225466198
$ws.Range("D8").Value = 0.30776620003126054
$ws.Range("E8").Value = 0.57151351471381817

$ws = $wb.Worksheets.Item("alpha_student_motivation")
$ws.Range("B2").Value = 0.94879129730325307
$ws.Range("B5").Value = 0.93251077638388924
$ws.Range("C5").Value = 0.87853374268790541
$ws.Range("D5").Value = 0.82029816548400292
$ws.Range("E5").Value = 0.93194649541525321
$ws.Range("B6").Value = 0.87021774315021971
$ws.Range("C6").Value = 0.77377575858429226
$ws.Range("D6").Value = 0.89889402367670834
$ws.Range("E6").Value = 0.96386217447067379
$ws.Range("B7").Value = 0.96807253793751136
$ws.Range("C7").Value = 0.94150521098901219
$ws.Range("D7").Value = 0.77628584993049388
$ws.Range("E7").Value = 0.91235734922482714
$ws.Range("B8").Value = 0.95472308681019069
$ws.Range("C8").Value = 0.9170570109632542
$ws.Range("D8").Value = 0.79430254063074579
$ws.Range("E8").Value = 0.92053733462154352
